# chore: update Sheets via scheduled runner
# Refresh cached Universalis market-price figures (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ columns, H:N) across the per-job leve tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 19763.334
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 19763.334
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = 19763.334
$ws.Range("N7").Value = -19987.334
$ws.Range("H14").Value = 19763.334
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 19763.334
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = ""
$ws.Range("M14").Value = 19763.334
$ws.Range("N14").Value = -20145.334
$ws.Range("H16").Value = 35249.75
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 35249.75
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = ""
$ws.Range("M16").Value = 35249.75
$ws.Range("N16").Value = -35709.75
$ws.Range("H126").Value = 46920
$ws.Range("J126").Value = 46920
$ws.Range("L126").Value = 46920
$ws.Range("N126").Value = -56800

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 4000.5
$ws.Range("J8").Value = 4000.5
$ws.Range("L8").Value = 4000.5
$ws.Range("N8").Value = -4288.5
$ws.Range("H10").Value = 72502.5
$ws.Range("J10").Value = 72502.5
$ws.Range("L10").Value = 72502.5
$ws.Range("N10").Value = -72842.5
$ws.Range("H12").Value = 47669.332
$ws.Range("I12").Value = 3000
$ws.Range("J12").Value = 70004
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 70004
$ws.Range("M12").Value = -2827
$ws.Range("N12").Value = -70350
$ws.Range("H104").Value = 42201
$ws.Range("J104").Value = 42201
$ws.Range("L104").Value = 42201
$ws.Range("N104").Value = -49189
$ws.Range("H117").Value = 46523.715
$ws.Range("J117").Value = 46523.715
$ws.Range("L117").Value = 46523.715
$ws.Range("N117").Value = -55701.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 47747
$ws.Range("J117").Value = 47747
$ws.Range("L117").Value = 47747
$ws.Range("N117").Value = -56925

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 17073.133
$ws.Range("I12").Value = 2700
$ws.Range("K12").Value = 2700
$ws.Range("M12").Value = -2530
$ws.Range("H19").Value = 699.5
$ws.Range("I19").Value = 99
$ws.Range("K19").Value = 99
$ws.Range("M19").Value = 71
$ws.Range("H24").Value = 699.5
$ws.Range("I24").Value = 99
$ws.Range("K24").Value = 99
$ws.Range("M24").Value = 71
$ws.Range("H104").Value = 29997.666
$ws.Range("J104").Value = 29997.666
$ws.Range("L104").Value = 29997.666
$ws.Range("N104").Value = -35239.666
$ws.Range("H109").Value = 27118.182
$ws.Range("J109").Value = 27118.182
$ws.Range("L109").Value = 27118.182
$ws.Range("N109").Value = -29198.182
$ws.Range("H115").Value = 28930.666
$ws.Range("J115").Value = 28930.666
$ws.Range("L115").Value = 28930.666
$ws.Range("N115").Value = -31280.666
$ws.Range("H116").Value = 47668.5
$ws.Range("J116").Value = 47668.5
$ws.Range("L116").Value = 47668.5
$ws.Range("N116").Value = -56846.5
$ws.Range("H120").Value = 32611.727
$ws.Range("J120").Value = 32611.727
$ws.Range("L120").Value = 32611.727
$ws.Range("N120").Value = -39869.727
$ws.Range("H141").Value = 7910.2856
$ws.Range("J141").Value = 7910.2856
$ws.Range("L141").Value = 7910.2856
$ws.Range("N141").Value = -18270.2856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3815.487
$ws.Range("J131").Value = 1639.2963
$ws.Range("L131").Value = 4917.8889
$ws.Range("N131").Value = -14997.8889

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 25447.5
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 25447.5
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = 25447.5
$ws.Range("N6").Value = -25673.5
$ws.Range("H9").Value = 1400
$ws.Range("I9").Value = 1400
$ws.Range("K9").Value = 1400
$ws.Range("M9").Value = -1230
$ws.Range("H16").Value = 25447.5
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 25447.5
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = ""
$ws.Range("M16").Value = 25447.5
$ws.Range("N16").Value = -25947.5
$ws.Range("H19").Value = 3781.6667
$ws.Range("I19").Value = 2195
$ws.Range("J19").Value = 4575
$ws.Range("K19").Value = 2195
$ws.Range("L19").Value = 4575
$ws.Range("M19").Value = -1907
$ws.Range("N19").Value = -5151
$ws.Range("H130").Value = 45993
$ws.Range("J130").Value = 45993
$ws.Range("L130").Value = 45993
$ws.Range("N130").Value = -56033

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 9700
$ws.Range("I74").Value = 9700
$ws.Range("K74").Value = 9700
$ws.Range("M74").Value = -8702
$ws.Range("H75").Value = 38900
$ws.Range("J75").Value = 38900
$ws.Range("L75").Value = 38900
$ws.Range("N75").Value = -40772
$ws.Range("H76").Value = 16491.666
$ws.Range("J76").Value = 19922.25
$ws.Range("L76").Value = 19922.25
$ws.Range("N76").Value = -20598.25
$ws.Range("H77").Value = 9700
$ws.Range("I77").Value = 9700
$ws.Range("K77").Value = 29100
$ws.Range("M77").Value = -24108
$ws.Range("H78").Value = 38900
$ws.Range("J78").Value = 38900
$ws.Range("L78").Value = 116700
$ws.Range("N78").Value = -126060
$ws.Range("H79").Value = 16491.666
$ws.Range("J79").Value = 19922.25
$ws.Range("L79").Value = 19922.25
$ws.Range("N79").Value = -22262.25
$ws.Range("H88").Value = 44185
$ws.Range("J88").Value = 44185
$ws.Range("L88").Value = 44185
$ws.Range("N88").Value = -45041
$ws.Range("H91").Value = 44185
$ws.Range("J91").Value = 44185
$ws.Range("L91").Value = 44185
$ws.Range("N91").Value = -47149
$ws.Range("H97").Value = 34996
$ws.Range("J97").Value = 34996
$ws.Range("L97").Value = 34996
$ws.Range("N97").Value = -36978
$ws.Range("H110").Value = 32357.334
$ws.Range("J110").Value = 32357.334
$ws.Range("L110").Value = 32357.334
$ws.Range("N110").Value = -40537.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 41755.4
$ws.Range("J16").Value = 41755.4
$ws.Range("L16").Value = 41755.4
$ws.Range("N16").Value = -42339.4
$ws.Range("H18").Value = 9534.333000000001
$ws.Range("J18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("N18").Value = -10346
$ws.Range("H119").Value = 48694
$ws.Range("J119").Value = 48694
$ws.Range("L119").Value = 48694
$ws.Range("N119").Value = -58370
$ws.Range("H122").Value = 1504805.6
$ws.Range("J122").Value = 958
$ws.Range("L122").Value = 2874
$ws.Range("N122").Value = -7774
$ws.Range("H123").Value = 31162.5
$ws.Range("J123").Value = 31162.5
$ws.Range("L123").Value = 31162.5
$ws.Range("N123").Value = -40962.5
